# SAP sales invoice generation process update:
# Remove the first three document numbers (old rows A2:A4) from the list,
# shifting all remaining SAP order numbers up by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the three obsolete order rows (mirrors selecting whole rows 2-4 in
# the UI before deleting them) and delete them, shifting the rest of the
# data up.
$ws.Range("A2:A4").EntireRow.Select()
$ws.Range("A2:A4").EntireRow.Delete()
